$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "43.879.43"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "2.258.27"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("E4").Value = "  -0.08%  "

Set-TextValue $ws.Range("D5") "230.41"
$ws.Range("E5").Value = "  -0.05%  "

Set-TextValue $ws.Range("D6") "0.643"
$ws.Range("E6").Value = "  +2.51%  "

Set-TextValue $ws.Range("D7") "64.34"
$ws.Range("E7").Value = "  +4.92%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  +7.12%  "

Set-TextValue $ws.Range("D10") "0.0982"
$ws.Range("E10").Value = "  +5.15%  "

Set-TextValue $ws.Range("D11") "57.00"
$ws.Range("E11").Value = "  -1.63%  "

Set-TextValue $ws.Range("D12") "27.03"
$ws.Range("E12").Value = "  +14.95%  "

Set-TextValue $ws.Range("D13") "0.105"
$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("D14").Value = "2.592.63"
$ws.Range("E14").Value = "  -0.35%  "

Set-TextValue $ws.Range("D15") "15.64"
$ws.Range("E15").Value = "  +0.41%  "

Set-TextValue $ws.Range("D16") "6.10"
$ws.Range("E16").Value = "  +5.34%  "

Set-TextValue $ws.Range("D17") "0.834"
$ws.Range("E17").Value = "  +3.26%  "

$ws.Range("D18").Value = "2.256.15"
$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("D19").Value = "43.739.57"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").Value = "0.0₃0998"
$ws.Range("E20").Value = "  +7.05%  "

$ws.Range("E21").Value = "  +0.69%  "

$ws.Range("E22").Value = "  -2.38%  "

Set-TextValue $ws.Range("D23") "251.32"
$ws.Range("E23").Value = "  -1.24%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("E25").Value = "  -4.38%  "

$ws.Range("E26").Value = "  -2.53%  "

$ws.Range("E27").Value = "  +2.44%  "

Set-TextValue $ws.Range("D29") "171.25"
$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("E30").Value = "  -0.72%  "

Set-TextValue $ws.Range("D31") "20.91"
$ws.Range("E31").Value = "  +2.10%  "

$ws.Range("E32").Value = "  -3.44%  "

$ws.Range("E33").Value = "  +2.93%  "

Set-TextValue $ws.Range("D34") "0.0707"
$ws.Range("E34").Value = "  +7.15%  "

Set-TextValue $ws.Range("D35") "4.78"
$ws.Range("E35").Value = "  +0.31%  "

Set-TextValue $ws.Range("D36") "4.90"
$ws.Range("E36").Value = "  -3.28%  "

Set-TextValue $ws.Range("D37") "3.76"
$ws.Range("E37").Value = "  +4.37%  "

$ws.Range("E38").Value = "  +0.05%  "

Set-TextValue $ws.Range("D39") "2.29"
$ws.Range("E39").Value = "  -4.04%  "

Set-TextValue $ws.Range("D40") "0.0260"
$ws.Range("E40").Value = "  +4.16%  "

$ws.Range("E41").Value = "  +0.04%  "

Set-TextValue $ws.Range("D42") "0.000225"
$ws.Range("E42").Value = "  -1.57%  "

Set-TextValue $ws.Range("D43") "0.0972"
$ws.Range("E43").Value = "  -1.06%  "

Set-TextValue $ws.Range("D44") "17.36"
$ws.Range("E44").Value = "  +4.48%  "

Set-TextValue $ws.Range("D45") "8.20"
$ws.Range("E45").Value = "  -6.00%  "

Set-TextValue $ws.Range("D46") "97.83"
$ws.Range("E46").Value = "  -0.27%  "

$ws.Range("E47").Value = "  -0.97%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "2.38"
$ws.Range("E48").Value = "  +5.92%  "

$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D49") "4.40"
$ws.Range("E49").Value = "  -2.92%  "

Set-TextValue $ws.Range("D50") "10.14"
$ws.Range("E50").Value = "  +6.07%  "

$ws.Range("D51").Value = "1.436.46"
$ws.Range("E51").Value = "  -2.42%  "
